# Applies the "network_diagram" edit:
#  - Adds 5 new "People" rows (63-67) under the OIT / Depends On / People
#    block, mirroring the existing block at rows 12-16, with 5 new staff
#    names and position/salary notes loaded into the right-hand container.
#  - Extends the formatted-but-empty column banding (A & E fills) down
#    through row 87 so the sheet keeps its look below the new data.
#  - Leaves the final selection on G81, matching the author's last click.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1) Build rows 63-67 by cloning the formatting of the existing,
#        identical "OIT / Depends On / People" block (rows 12-16) ---
$ws.Range("A12:G12").Copy()
$ws.Range("A63:G67").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A63:G67").RowHeight = 42.75

$peopleRows = @(63, 64, 65, 66, 67)
$names = @("Liam Anderson", "Emma Roberts", "Noah Carter", "Olivia Mitchell", "William Turner")
$notes = @(
  "POS #: Staff Aug`nTitle:  DB Admin - Lead`nS+B: `$130,750",
  "POS #: Staff Aug`nTitle:  DB Admin - Lead`nS+B: `$130,751",
  "POS #: Staff Aug`nTitle:  DB Admin - Lead`nS+B: `$130,752",
  "POS #: Staff Aug`nTitle:  DB Admin - Lead`nS+B: `$130,753",
  "POS #: Staff Aug`nTitle:  DB Admin - Lead`nS+B: `$130,754"
)
$orgDescr = "Office of Information Technology (OIT)`nThe department managing the financial responsibilities for the State of Florida. "

foreach ($r in $peopleRows) {
    $ws.Range("A$r").Value = "Organization"
    $ws.Range("B$r").Value = "OIT"
    $ws.Range("C$r").Value = $orgDescr
    $ws.Range("D$r").Value = "Depends On"
    $ws.Range("E$r").Value = "People"
}
for ($i = 0; $i -lt $peopleRows.Length; $i++) {
    $ws.Range("F" + $peopleRows[$i]).Value = $names[$i]
}
for ($i = 0; $i -lt $peopleRows.Length; $i++) {
    $ws.Range("G" + $peopleRows[$i]).Value = $notes[$i]
}

# --- 2) Extend the column banding (A & E fills) as empty formatted rows
#        68-87, mirroring the formatting already used in the block above ---
$ws.Range("A12").Copy()
$ws.Range("A68:A87").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E12").Copy()
$ws.Range("E68:E87").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# --- 3) Restore a sane view: scroll near the bottom of the frozen pane
#        and leave the final selection where the author left it ---
$excel.ActiveWindow.ScrollRow = 53
$ws.Range("G81").Select() | Out-Null
